# #5: insurance, claim, debt, investment done
#
# The "保險" (insurance) sheet used to have an ad-hoc 5-column layout
# (company / product / owner / ??? / free-text note). Bring it in line
# with the other property sheets: B=company name, C=product name,
# D=owner, E=property_category, then the standard trailing metadata
# columns (category, date, legislator_name, legislator_id, source_file,
# index).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("保險")

# --- Header row (row 1) ------------------------------------------------
# B1 ("company") is already correct; C1/D1/E1 are repointed to the new
# generic column names, and F1:K1 are brand-new header cells that need
# the same bold/centered/bordered look as the rest of the header.
$ws.Cells.Item(1, 3).Value = "name"
$ws.Cells.Item(1, 4).Value = "owner"
$ws.Cells.Item(1, 5).Value = "property_category"

$headerRange = $ws.Range("F1:K1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$ws.Cells.Item(1, 6).Value = "category"
$ws.Cells.Item(1, 7).Value = "date"
$ws.Cells.Item(1, 8).Value = "legislator_name"
$ws.Cells.Item(1, 9).Value = "legislator_id"
$ws.Cells.Item(1, 10).Value = "source_file"
$ws.Cells.Item(1, 11).Value = "index"

# --- Row 2 (index 87) ---------------------------------------------------
$ws.Cells.Item(2, 2).Value = "中國人壽"
$ws.Cells.Item(2, 3).Value = "喜悅人生變額壽險"
$ws.Cells.Item(2, 4).Value = "田秋堇"
$ws.Cells.Item(2, 5).Value = "insurance"
$ws.Cells.Item(2, 6).Value = "normal"
$ws.Cells.Item(2, 7).NumberFormat = "@"
$ws.Cells.Item(2, 7).Value = "2012-04-10"
$ws.Cells.Item(2, 8).Value = "田秋堇"
$ws.Cells.Item(2, 9).Value = 1316
$ws.Cells.Item(2, 10).Value = "tmp9b251"
$ws.Cells.Item(2, 11).Value = 87

# --- Row 3 (index 88) ---------------------------------------------------
$ws.Cells.Item(3, 2).Value = "中國人壽"
$ws.Cells.Item(3, 3).Value = "喜悅人生變額壽險"
$ws.Cells.Item(3, 4).Value = "劉守成"
$ws.Cells.Item(3, 5).Value = "insurance"
$ws.Cells.Item(3, 6).Value = "normal"
$ws.Cells.Item(3, 7).NumberFormat = "@"
$ws.Cells.Item(3, 7).Value = "2012-04-10"
$ws.Cells.Item(3, 8).Value = "田秋堇"
$ws.Cells.Item(3, 9).Value = 1316
$ws.Cells.Item(3, 10).Value = "tmp9b251"
$ws.Cells.Item(3, 11).Value = 88

# --- Row 4 (index 89) ---------------------------------------------------
$ws.Cells.Item(4, 2).Value = "富邦人壽"
$ws.Cells.Item(4, 3).Value = "歡喜年年終身乙型"
$ws.Cells.Item(4, 4).Value = "劉守成"
$ws.Cells.Item(4, 5).Value = "insurance"
$ws.Cells.Item(4, 6).Value = "normal"
$ws.Cells.Item(4, 7).NumberFormat = "@"
$ws.Cells.Item(4, 7).Value = "2012-04-10"
$ws.Cells.Item(4, 8).Value = "田秋堇"
$ws.Cells.Item(4, 9).Value = 1316
$ws.Cells.Item(4, 10).Value = "tmp9b251"
$ws.Cells.Item(4, 11).Value = 89
